# Update "想去人数" (want-to-go count) figures for two events whose rows
# appear on both the "展览" sheet and the aggregate "全部类型" sheet.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 6593
    $ws.Range("F6").Value = 131
}
